$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell text: userActions -> apiKeyActions
$ws.Range("A1").Value = "button_apiKeyActions_trNthChild"
$ws.Range("B1").Value = "button_apiKeyActions_trNthChild_1"
$ws.Range("C1").Value = "button_apiKeyActions_trNthChild_2"

# Update data path value
$ws.Range("D2").Value = "Data Files/AI-Generated/Common/createApiKeyForAccess-test-data"

# Update column widths
# Note: Excel's ColumnWidth property stores the value with an added ~0.8333
# padding offset when round-tripped through the OOXML "width" attribute, so
# we subtract that offset here to land exactly on the target stored widths
# (33, 35, 35, 64).
$ws.Columns.Item(1).ColumnWidth = 32.166666666666664
$ws.Columns.Item(2).ColumnWidth = 34.166666666666664
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664
$ws.Columns.Item(4).ColumnWidth = 63.16666666666667
